# edit.ps1
# Applies the "rerrange the script to be more lean" commit to Data_analysis.xlsx
# (workbook already open as $excel.ActiveWorkbook)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Branchmark Airplane")

# ---------------------------------------------------------------------------
# 1. Add a new blank worksheet "Sheet2" as the very last sheet in the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count())
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet2"
$newSheet.Range("D9").Select()

# Re-activate the original sheet so it stays the "tabSelected" sheet
$ws.Activate()

# ---------------------------------------------------------------------------
# 2. Update the benchmark table (B6:H18): add left-alignment to every cell
#    that already carries the thin-border style, keeping borders / number
#    formats / wrap / formulas untouched.
# ---------------------------------------------------------------------------
$ws.Range("B6:H18").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# 3. Refresh "Fuel saved" row (row 16) with the re-run simulation numbers
# ---------------------------------------------------------------------------
$ws.Range("C16").Value() = 16353.2425826376
$ws.Range("D16").Value() = 49214.444181069703
$ws.Range("E16").Value() = 66160.043385847996
$ws.Range("F16").Value() = 73861.227986158905
$ws.Range("G16").Value() = 173694.952401569
$ws.Range("H16").Value() = 206932.244301542

# ---------------------------------------------------------------------------
# 4. Wipe the old scratch / scribble cells that used to live below the table
#    (J20:L26 helper values, and the stray F27:F32 column) - the sheet is
#    being "releaned" and re-purposed for the new data below.
# ---------------------------------------------------------------------------
$ws.Range("J20:L26").ClearContents()
$ws.Range("F27:F32").ClearContents()

# ---------------------------------------------------------------------------
# 5. New row 20: "Take-off distance[ft] 1.1 margin / (Simulations)"
# ---------------------------------------------------------------------------
$ws.Range("B16").Copy()
$ws.Range("B20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B20").Value() = "Take-off distance[ft] 1.1 margin`n(Simulations)"
$ws.Rows.Item(20).RowHeight = 42.75

$ws.Range("C20").Value() = 3433.5367253872801
$ws.Range("D20").Value() = 7762.6418483978396
$ws.Range("E20").Value() = 8659.8757447735898
$ws.Range("F20").Value() = 7626.3734160430404
$ws.Range("G20").Value() = 5679.3101547331999
$ws.Range("H20").Value() = 6704.4421894425705

# ---------------------------------------------------------------------------
# 6. "Optimization using simulated annealing" header (merged B24:C24)
# ---------------------------------------------------------------------------
$ws.Range("B24:C24").Merge()
$ws.Range("B24").Value() = "Optimization using simulated annealing"
$ws.Range("B24:C24").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 7. Header row for the new optimization table (row 25) - wrap text, no
#    border, custom row formatting
# ---------------------------------------------------------------------------
$ws.Range("B25").Value() = "Case #"
$ws.Range("C25").Value() = "Range(ship)[km]"
$ws.Range("D25").Value() = "Range(island)[km]"
$ws.Range("E25").Value() = "Run way limits[ft]"
$ws.Range("F25").Value() = "Ship runway limits[ft]"
$ws.Range("G25").Value() = "Wing span[ft]"
$ws.Range("H25").Value() = "thr to weight ratio[-]"
$ws.Range("I25").Value() = "aspect ratio[-]"
$ws.Range("J25").Value() = "sweep angle[deg]"
$ws.Range("K25").Value() = "max takeoff weight[lb]"
$ws.Range("L25").Value() = "FUEL SAVED"

$ws.Range("B25:L25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 28.5

# ---------------------------------------------------------------------------
# 8. Optimization results table (rows 26-31)
# ---------------------------------------------------------------------------
$ws.Range("B26").Value() = 1
$ws.Range("C26").Value() = 500
$ws.Range("D26").Value() = 900
$ws.Range("E26").Value() = 9000
$ws.Range("F26").Value() = 205
$ws.Range("G26").Value() = 120.96161081621101
$ws.Range("H26").Value() = 0.35
$ws.Range("I26").Value() = 5.57270323280882
$ws.Range("J26").Value() = 23.522281247699901
$ws.Range("K26").Value() = 137324.62388941899

$ws.Range("B27").Value() = 2
$ws.Range("C27").Value() = 500
$ws.Range("D27").Value() = 900
$ws.Range("E27").Value() = 9000
$ws.Range("F27").Value() = 245
$ws.Range("G27").Value() = 105.556417897317
$ws.Range("H27").Value() = 0.34982551871945899
$ws.Range("I27").Value() = 5.6132254754894104
$ws.Range("J27").Value() = 0.31863630739584903
$ws.Range("K27").Value() = 118574.770218673

$ws.Range("B28").Value() = 3
$ws.Range("C28").Value() = 500
$ws.Range("D28").Value() = 900
$ws.Range("E28").Value() = 9000
$ws.Range("F28").Value() = 285
$ws.Range("G28").Value() = 127.26848721229
$ws.Range("H28").Value() = 0.35
$ws.Range("I28").Value() = 5.5428819204567104
$ws.Range("J28").Value() = 21.928792935022699
$ws.Range("K28").Value() = 153024.91898027199

$ws.Range("B29").Value() = 4
$ws.Range("C29").Value() = 500
$ws.Range("D29").Value() = 900
$ws.Range("E29").Value() = 9000
$ws.Range("F29").Value() = 330
$ws.Range("G29").Value() = 131.82494344988299
$ws.Range("H29").Value() = 0.35
$ws.Range("I29").Value() = 5.9679445831296798
$ws.Range("J29").Value() = 24.371864113890702
$ws.Range("K29").Value() = 157925.982587951

$ws.Range("B30").Value() = 5
$ws.Range("C30").Value() = 500
$ws.Range("D30").Value() = 900
$ws.Range("E30").Value() = 9000
$ws.Range("F30").Value() = 415
$ws.Range("G30").Value() = 128.2603542905
$ws.Range("H30").Value() = 0.35
$ws.Range("I30").Value() = 5.57709068728267
$ws.Range("J30").Value() = 19.585818931124098
$ws.Range("K30").Value() = 169949.95488569399

$ws.Range("B31").Value() = 6
$ws.Range("C31").Value() = 500
$ws.Range("D31").Value() = 900
$ws.Range("E31").Value() = 9000
$ws.Range("F31").Value() = 9000
$ws.Range("G31").Value() = 128.2603542905
$ws.Range("H31").Value() = 0.35
$ws.Range("I31").Value() = 10.1
$ws.Range("J31").Value() = 19.585818931124098
$ws.Range("K31").Value() = 169949.95488569399

# ---------------------------------------------------------------------------
# 9. Leftover helper values (rows 48-52, column D)
# ---------------------------------------------------------------------------
$ws.Range("D48").Value() = 7762.6418483978396
$ws.Range("D49").Value() = 8659.8757447735898
$ws.Range("D50").Value() = 7626.3734160430404
$ws.Range("D51").Value() = 5679.3101547331999
$ws.Range("D52").Value() = 6704.4421894425705

# ---------------------------------------------------------------------------
# 10. Column width re-tune (widths widened to fit new content)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(7).ColumnWidth = 11
$ws.Columns.Item(8).ColumnWidth = 16
$ws.Columns.Item(9).ColumnWidth = 11
$ws.Columns.Item(10).ColumnWidth = 14
$ws.Columns.Item(11).ColumnWidth = 18

# ---------------------------------------------------------------------------
# 11. Final selection, matching the author's last cursor position
# ---------------------------------------------------------------------------
$ws.Range("B21").Select()
